# Update data values per diff ("Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.413
$ws.Range("E4").Value = 16.379
$ws.Range("C6").Value = -12.661
$ws.Range("B7").Value = 5.865
$ws.Range("D7").Value = -7.968999999999999
$ws.Range("A9").Value = -21.881
$ws.Range("D10").Value = -8.415000000000001
$ws.Range("B12").Value = 5.782999999999999
$ws.Range("D13").Value = -7.833
$ws.Range("E13").Value = 16.616
$ws.Range("B14").Value = 6.229
$ws.Range("C15").Value = -12.909
$ws.Range("D16").Value = -8.657
$ws.Range("A18").Value = -21.993
$ws.Range("A20").Value = -20.931
$ws.Range("D20").Value = -7.874
$ws.Range("D24").Value = -7.547
$ws.Range("B26").Value = 5.879
$ws.Range("A27").Value = -21.611
$ws.Range("B27").Value = 6.446000000000001
$ws.Range("E27").Value = 16.737
$ws.Range("B29").Value = 5.566
$ws.Range("E29").Value = 17.108
$ws.Range("C33").Value = -11.314
$ws.Range("A35").Value = -20.569
$ws.Range("C35").Value = -12.357
$ws.Range("E35").Value = 16.204
$ws.Range("B37").Value = 8.270000000000001
$ws.Range("B38").Value = 4.83
$ws.Range("C38").Value = -12.595
$ws.Range("D39").Value = -7.634
$ws.Range("E40").Value = 16.618
$ws.Range("C43").Value = -12.494
$ws.Range("C44").Value = -12.008
$ws.Range("C47").Value = -11.418
$ws.Range("D47").Value = -7.192
$ws.Range("D48").Value = -7.228999999999999
$ws.Range("B51").Value = 5.824
$ws.Range("C51").Value = -12.418
$ws.Range("B52").Value = 5.478
$ws.Range("D52").Value = -7.668000000000001
$ws.Range("B55").Value = 6.403
$ws.Range("D56").Value = -7.825
$ws.Range("C57").Value = -12.742
$ws.Range("E57").Value = 16.844
$ws.Range("C63").Value = -11.796
$ws.Range("A69").Value = -21.701
$ws.Range("B69").Value = 6.616
$ws.Range("B70").Value = 6.114
$ws.Range("C70").Value = -11.567
$ws.Range("A76").Value = -20.738
$ws.Range("A78").Value = -20.609
$ws.Range("B81").Value = 6.154
$ws.Range("A82").Value = -21.837
$ws.Range("A83").Value = -20.61
$ws.Range("B83").Value = 7.384
$ws.Range("D84").Value = -8.294
$ws.Range("E85").Value = 16.485
$ws.Range("C88").Value = -12.583
$ws.Range("A93").Value = -22.021
$ws.Range("C99").Value = -12.399
$ws.Range("D100").Value = -8.311000000000002
$ws.Range("D101").Value = -7.831
$ws.Range("B102").Value = 7.402999999999999
